$p = $ppt.ActivePresentation

# 1) Update the "fixed" date/time footer placeholder from 22/08/2017 to
#    25/08/2017 on the slide master and on every slide layout (each one
#    carries its own copy of the Date Placeholder shape).
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "25/08/2017"
    }
}

for ($k = 1; $k -le $m.CustomLayouts.Count; $k++) {
    $cl = $m.CustomLayouts.Item($k)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "25/08/2017"
        }
    }
}

# 2) Update the GitHub repository link shown on the "Link do GitHub" slide.
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -like "*github.com*") {
        $sh.TextFrame.TextRange.Text = " https://github.com/betofr1/co-design.git"
    }
}
